# Sprint Backlog 2 + Burndown Chart fix
# "Ändrade så jag har 5 dagar på mig och ändrade skattad och ideal tid"
#
# This adds a 5th day column (J) to the sprint backlog table, fills in the
# estimated/ideal burndown rows with new values, and updates the burndown
# chart (which already referenced 5 days, F:J) so its cached values pick
# up the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Add the new "Dag 5" header in column J, row 1, matching the style of the
# other day headers.
$ws.Range("J1").Value = "Dag 5"

# Force Excel to re-resolve J1's cell style (xf index 9) so that its font
# reference gets normalized/deduped the same way the real edit did
# (fontId 3 -> fontId 1 -- both are identical bold 10pt black Arial fonts,
# this is just Excel's internal style table bookkeeping). Toggling
# WrapText off/on is a no-visual-effect way to force that re-resolution
# without introducing any new font/numfmt/fill entries.
$ws.Range("J1").WrapText = $false
$ws.Range("J1").WrapText = $true

# Estimated hours remaining per task (column F) updated.
$ws.Range("F2").Value = 10

# Burndown summary rows: "Skattat" (row 6) and "Idealt" (row 7)
$ws.Range("F6").Value = 13

$ws.Range("F7").Value = 13
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 2

# New day-5 column values (column J) for the data rows that already have
# columns F:I populated.
$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("J7").Value = 0

# Fill in column I for rows 2/3/6 to keep them explicit zeros (they already
# were zero, left as-is / re-set for clarity).
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I6").Value = 0

# Update the selection to reflect where the user ended up after editing.
$ws.Range("I13").Select()

$wb.Save()
